$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.548.08"
$ws.Range("E2").Value = "  +0.60%  "
$ws.Range("D3").Value = "2.977.57"
$ws.Range("E3").Value = "  +1.97%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").Value = "'381.82"
$ws.Range("E5").Value = "  +2.31%  "
$ws.Range("D6").Value = "'103.88"
$ws.Range("E6").Value = "  +0.33%  "
$ws.Range("E7").Value = "  +0.70%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D9").Value = "'0.592"
$ws.Range("E9").Value = "  +0.29%  "
$ws.Range("D10").Value = "'37.03"
$ws.Range("E10").Value = "  -0.06%  "
$ws.Range("E11").Value = "  +0.03%  "
$ws.Range("D12").Value = "'0.0846"
$ws.Range("E12").Value = "  +1.03%  "
$ws.Range("D13").Value = "3.446.53"
$ws.Range("E13").Value = "  +1.99%  "
$ws.Range("D14").Value = "'18.27"
$ws.Range("E14").Value = "  -0.36%  "
$ws.Range("D15").Value = "'7.59"
$ws.Range("E15").Value = "  +2.39%  "
$ws.Range("D16").Value = "2.978.17"
$ws.Range("E16").Value = "  +1.98%  "
$ws.Range("D17").Value = "'0.992"
$ws.Range("E17").Value = "  +6.62%  "
$ws.Range("D18").Value = "51.485.00"
$ws.Range("E19").Value = "  -1.52%  "
$ws.Range("D20").Value = "'7.41"
$ws.Range("E20").Value = "  +1.69%  "
$ws.Range("D21").Value = "'12.79"
$ws.Range("E21").Value = "  -1.55%  "
$ws.Range("D22").Value = "0.0₃0963"
$ws.Range("E22").Value = "  +1.65%  "
$ws.Range("D23").Value = "'69.08"
$ws.Range("E23").Value = "  +0.84%  "
$ws.Range("D24").Value = "'261.71"
$ws.Range("E24").Value = "  +0.39%  "
$ws.Range("D25").Value = "'2.91"
$ws.Range("E25").Value = "  +7.95%  "
$ws.Range("D26").Value = "'8.22"
$ws.Range("E26").Value = "  +14.82%  "
$ws.Range("D27").Value = "'7.69"
$ws.Range("E27").Value = "  +18.15%  "
$ws.Range("D28").Value = "'0.115"
$ws.Range("E28").Value = "  +12.71%  "
$ws.Range("D29").Value = "'0.169"
$ws.Range("E29").Value = "  -1.06%  "
$ws.Range("B30").Value = "LEO"
$ws.Range("C30").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D30").Value = "'4.13"
$ws.Range("E30").Value = "  +0.71%  "
$ws.Range("B31").Value = "Dai"
$ws.Range("C31").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D31").Value = "'1.00"
$ws.Range("E31").Value = "  -0.04%  "
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").Value = "'25.95"
$ws.Range("E32").Value = "  +0.55%  "
$ws.Range("B33").Value = "Cosmos"
$ws.Range("C33").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D33").Value = "'9.86"
$ws.Range("E33").Value = "  -0.35%  "
$ws.Range("B34").Value = "InjectiveProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D34").Value = "'34.53"
$ws.Range("E34").Value = "  -0.28%  "
$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D35").Value = "'50.94"
$ws.Range("E35").Value = "  -0.53%  "
$ws.Range("B36").Value = "Toncoin"
$ws.Range("C36").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D36").Value = "'2.06"
$ws.Range("E36").Value = "  -2.38%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "'0.0449"
$ws.Range("E37").Value = "  +5.75%  "
$ws.Range("B38").Value = "FirstDigitalUSD"
$ws.Range("C38").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D38").Value = "'1.00"
$ws.Range("E38").Value = "  -0.14%  "
$ws.Range("B39").Value = "LidoDAOToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D39").Value = "'3.00"
$ws.Range("E39").Value = "  +0.28%  "
$ws.Range("B40").Value = "Celestia"
$ws.Range("C40").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D40").Value = "'16.94"
$ws.Range("E40").Value = "  -0.86%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "'2.57"
$ws.Range("E41").Value = "  -0.62%  "
$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").Value = "'0.116"
$ws.Range("E42").Value = "  +2.35%  "
$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").Value = "'1.83"
$ws.Range("E43").Value = "  -1.04%  "
$ws.Range("B44").Value = "Monero"
$ws.Range("C44").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D44").Value = "'122.32"
$ws.Range("E44").Value = "  +2.28%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "'21.69"
$ws.Range("E45").Value = "  -1.37%  "
$ws.Range("B46").Value = "TheGraph"
$ws.Range("C46").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D46").Value = "'0.276"
$ws.Range("E46").Value = "  +12.33%  "
$ws.Range("B47").Value = "WEMIXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D47").Value = "'2.05"
$ws.Range("E47").Value = "  -2.17%  "
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").Value = "'3.31"
$ws.Range("E49").Value = "  +4.26%  "
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").Value = "2.028.98"
$ws.Range("E50").Value = "  +0.21%  "
$ws.Range("B51").Value = "BEAM"
$ws.Range("C51").Value = "https://coinranking.com/coin/cYYMfXF4u+beam-beam"
$ws.Range("D51").Value = "'0.0331"
$ws.Range("E51").Value = "  +4.98%  "
